{"js": "// Insert a new run containing a single backtick (\"`\") character as the very\n// first run of the document's first paragraph (\"CS 3305: Data Structures\"),\n// matching the formatting (Arial, 14pt / 28 half-points, incl. complex-script\n// font + size) already used by the surrounding runs in that paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The target run is inserted at the very start of the title paragraph\n// (\"CS 3305: Data Structures\"), which is the first paragraph of the body.\n// Prefer locating it by its text so the script is resilient even if\n// paragraphs were ever reordered; fall back to the first paragraph.\nlet targetParagraph = paragraphs.items.find((p) => p.text.indexOf(\"CS 3305\") === 0);\nif (!targetParagraph) {\n  targetParagraph = paragraphs.items[0];\n}\nconst insertionPoint = targetParagraph.getRange(\"Start\");\n\n// Build a minimal OOXML package fragment describing the run to insert so we\n// control every run-property (ascii/hAnsi/cs font + sz/szCs) exactly,\n// matching the formatting already applied to the paragraph's other runs.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n                <w:sz w:val=\"28\"/>\n                <w:szCs w:val=\"28\"/>\n              </w:rPr>\n              <w:t>\\u0060</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionPoint.insertOoxml(ooxml, Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Insert a new run containing a single backtick (\"`\") character as the very\n# first run of the document's first paragraph (\"CS 3305: Data Structures\"),\n# matching the formatting (Arial, 14pt / 28 half-points, incl. complex-script\n# font + size) already used by the surrounding runs in that paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the very start of the \"CS 3305\" title so the edit does not depend\n# on the paragraph's ordinal position within the document.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"CS 3305\")\n\nif (-not $found) {\n    # Fallback: use the first paragraph of the document body.\n    $rng = $d.Paragraphs.Item(1).Range\n}\n\n$rng.Collapse(1)  # wdCollapseStart\n\n# Use InsertXML (raw WordprocessingML via the xmlPackage wrapper) so the\n# inserted run carries the exact run formatting (ascii/hAnsi/cs fonts and\n# sz/szCs sizes) rather than whatever the Font object partially applies.\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>`</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$rng.InsertXML($xml)\n"}
